$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
try {
    $ip = $tr.InsertAfter("")
    Write-Output "InsertAfter ok: [$($ip.Text)] len=$($ip.Length)"
} catch {
    Write-Output "InsertAfter failed: $_"
}
